$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header cell B2 from "unnamed: 1_level_1" to "total"
$ws.Range("B2").Value = "total"

# Remove the two blank "section header" rows (old row 8 then old row 5, highest index
# first so the remaining row number doesn't shift before we act on it).
$ws.Rows.Item(8).Delete()
$ws.Rows.Item(5).Delete()
